# Updates cryptos list values (price/volume) and swaps the Elrond/Aptos rows,
# matching the upstream data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.652.69'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '1.887.84'
$ws.Range("E3").Value = '  +1.00%  '
$ws.Range("D4").Formula = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Formula = "'247.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Formula = "'0.4729"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Formula = "'0.2918"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").Formula = "'0.06518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.60%  '
$ws.Range("D10").Formula = "'22.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.54%  '
$ws.Range("D11").Formula = "'0.07781"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").Value = '1.891.14'
$ws.Range("E12").Value = '  +1.14%  '
$ws.Range("D13").Formula = "'0.7400"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("D14").Formula = "'96.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("D15").Formula = "'5.232"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.72%  '
$ws.Range("D16").Formula = "'284.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.05%  '
$ws.Range("D17").Value = '30.629.19'
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("D18").Formula = "'13.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("D19").Formula = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").Formula = "'0.000007506"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '2.141.23'
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("D22").Formula = "'5.301"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.23%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Formula = "'6.246"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.25%  '
$ws.Range("D25").Formula = "'9.213"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.62%  '
$ws.Range("D26").Formula = "'164.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.63%  '
$ws.Range("D27").Formula = "'18.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("D28").Formula = "'1.914"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").Formula = "'1.344"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.44%  '
$ws.Range("D30").Formula = "'0.09759"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.40%  '
$ws.Range("D32").Formula = "'4.296"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("D33").Formula = "'4.179"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.56%  '
$ws.Range("D34").Formula = "'0.04899"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.91%  '
$ws.Range("D35").Formula = "'1.127"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.01%  '
$ws.Range("D36").Formula = "'0.6976"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.58%  '
$ws.Range("D37").Formula = "'2.716"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("D38").Formula = "'0.01895"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("D39").Formula = "'2.832"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.28%  '
$ws.Range("D40").Formula = "'76.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.60%  '
$ws.Range("E41").Value = '  +0.33%  '
$ws.Range("E42").Value = '  +1.80%  '
$ws.Range("D43").Formula = "'0.4274"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.19%  '
$ws.Range("D44").Formula = "'1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").Formula = "'0.8327"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.37%  '
$ws.Range("D46").Formula = "'101.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").Formula = "'9.597"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.99%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Formula = "'7.014"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.76%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Formula = "'35.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").Formula = "'904.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.77%  '
$ws.Range("D51").Formula = "'0.05774"
$ws.Range("D51").Style = "Normal"
